$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -7
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -5
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = -8
$ws.Range("F12").Value = -7
$ws.Range("F13").Value = 1
